# "Logged Week 17 data and fixed Simulate_Season.py tiebreaking method"
#
# This workbook stores two kinds of Week-by-week logs:
#   1) Space-separated running logs of per-game values on the YDS and ST
#      sheets (one new number/tuple gets appended to the string each week).
#   2) Season-to-date totals on the OFF / DEF / ST / TURNS / PEN sheets
#      (the box-score row for each team gets bumped by this week's numbers).
#
# Week 17 numbers are appended / added below.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# YDS sheet — append this week's Rushing/Passing yards-per-game tuple
# ------------------------------------------------------------------
$ydsWs = $wb.Worksheets.Item("YDS")
$ydsWs.Range("B2").Value = $ydsWs.Range("B2").Value2 + " 16 0 1 0 1 14 11 1 0 6 1 8 1 6 5 3"
$ydsWs.Range("C2").Value = $ydsWs.Range("C2").Value2 + " 0 4 5 -3 3 14 0 4 4 5 4 -1 6 2 5 21 2 7 5 3 2 4 3 1 0 12 0 9 35 -1 2 4 0 -3 0 0 -2 1 39"
$ydsWs.Range("B3").Value = $ydsWs.Range("B3").Value2 + " 20 0 8 14 14 11 15 21 0 10 6 6 45 1 2 9 8"
$ydsWs.Range("C3").Value = $ydsWs.Range("C3").Value2 + " 25 2 4 1 5 5 16 3 7 52 6 3 15"

# ------------------------------------------------------------------
# OFF sheet — add Week 17 totals onto the running season totals
# ------------------------------------------------------------------
$offWs = $wb.Worksheets.Item("OFF")
$offWs.Range("B2").Value = 2
$offWs.Range("C2").Value = 180
$offWs.Range("F2").Value = 66
$offWs.Range("G2").Value = 50
$offWs.Range("I2").Value = 8
$offWs.Range("J2").Value = 29
$offWs.Range("N2").Value = 16
$offWs.Range("O2").Value = 21
$offWs.Range("P2").Value = 12
$offWs.Range("C3").Value = 189
$offWs.Range("E3").Value = 27
$offWs.Range("F3").Value = 115
$offWs.Range("H3").Value = 28
$offWs.Range("I3").Value = 61
$offWs.Range("J3").Value = 69
$offWs.Range("L3").Value = 298
$offWs.Range("M3").Value = 193
$offWs.Range("Q3").Value = 496

# ------------------------------------------------------------------
# DEF sheet — add Week 17 totals onto the running season totals
# ------------------------------------------------------------------
$defWs = $wb.Worksheets.Item("DEF")
$defWs.Range("C2").Value = 190
$defWs.Range("E2").Value = 13
$defWs.Range("F2").Value = 65
$defWs.Range("G2").Value = 29
$defWs.Range("H2").Value = 9
$defWs.Range("J2").Value = 28
$defWs.Range("N2").Value = 19
$defWs.Range("O2").Value = 18
$defWs.Range("P2").Value = 11
$defWs.Range("B3").Value = 11
$defWs.Range("C3").Value = 183
$defWs.Range("D3").Value = 4
$defWs.Range("E3").Value = 49
$defWs.Range("F3").Value = 117
$defWs.Range("G3").Value = 27
$defWs.Range("I3").Value = 69
$defWs.Range("J3").Value = 49
$defWs.Range("L3").Value = 285
$defWs.Range("M3").Value = 189
$defWs.Range("Q3").Value = 511

# ------------------------------------------------------------------
# ST sheet — append this week's KO/PT per-game tuples, and bump the
# season totals for KO/PT attempts
# ------------------------------------------------------------------
$stWs = $wb.Worksheets.Item("ST")
$stWs.Range("B2").Value = 72
$stWs.Range("D2").Value = 73
$stWs.Range("J2").Value = 56
$stWs.Range("K2").Value = 53
$stWs.Range("N2").Value = 21
$stWs.Range("B3").Value = 51
$stWs.Range("B4").Value = $stWs.Range("B4").Value2 + " 64"
$stWs.Range("B5").Value = $stWs.Range("B5").Value2 + " 24"
$stWs.Range("B6").Value = $stWs.Range("B6").Value2 + " 16 18"
$stWs.Range("D3").Value = $stWs.Range("D3").Value2 + " 37 36 42 32"
$stWs.Range("D4").Value = $stWs.Range("D4").Value2 + " 9 8 11 0"
$stWs.Range("D5").Value = $stWs.Range("D5").Value2 + " 15 0 0"

# ------------------------------------------------------------------
# TURNS sheet — add Week 17 turnover totals
# ------------------------------------------------------------------
$turnsWs = $wb.Worksheets.Item("TURNS")
$turnsWs.Range("B3").Value = 7
$turnsWs.Range("D3").Value = 13
$turnsWs.Range("E3").Value = 8

# ------------------------------------------------------------------
# PEN sheet — add Week 17 "Intentional grounding" penalty count
# ------------------------------------------------------------------
$penWs = $wb.Worksheets.Item("PEN")
$penWs.Range("D4").Value = 5
